$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133; existing rows 133:194 shift down to 134:195.
$ws.Rows("133:133").Insert()

# Populate the newly inserted row 133 with the new weekly record.
$ws.Range("A133").Value = 5
$ws.Range("B133").Value = "Macroferia Regional de Talca"
$ws.Range("C133").Value = "Maule"
$ws.Range("D133").Value = 44572
$ws.Range("D133").NumberFormat = $ws.Range("D134").NumberFormat
$ws.Range("E133").Value = 7
$ws.Range("F133").Value = 100112008
$ws.Range("G133").Value = "Coliflor"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 2000
$ws.Range("K133").Value = 800
$ws.Range("L133").Value = 800
$ws.Range("M133").Value = 800
$ws.Range("N133").Value = "$/unidad"
$ws.Range("O133").Value = "Región del Maule"
$ws.Range("P133").Value = 800
$ws.Range("Q133").Value = 1
$ws.Range("R133").Value = "Hortaliza"
